# "ID and NRS pathology"
#
# 1. The first "Body" paragraph switches from the BodyText style to the
#    (BodyText-based) FirstParagraph style.
# 2. The APA template's heading / body-text / bibliography styles move
#    from the "CMU Serif*" family to Georgia.

$d = $word.ActiveDocument

# --- 1. Re-style the opening body paragraph -------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Body" -and $p.Style.NameLocal -eq "Body Text") {
        $p.Style = "FirstParagraph"
    }
}

# --- 2. Swap the template's serif font for Georgia -------------------------
$d.Styles("Heading1").Font.NameAscii = "Georgia"
$d.Styles("Heading1").Font.NameOther = "Georgia"

$d.Styles("Heading2").Font.NameAscii = "Georgia"
$d.Styles("Heading2").Font.NameOther = "Georgia"

$d.Styles("Heading3").Font.NameAscii = "Georgia"
$d.Styles("Heading3").Font.NameOther = "Georgia"

$d.Styles("BodyText").Font.NameAscii = "Georgia"
$d.Styles("BodyText").Font.NameOther = "Georgia"
$d.Styles("BodyText").Font.NameBi = "Times New Roman (Body CS)"
